$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.140870094299316
$ws.Range("B1").Value = 1.685660600662231
$ws.Range("C1").Value = 3.589289665222168
$ws.Range("D1").Value = 3.459939956665039
$ws.Range("E1").Value = 0.9598383903503418
